# Auto-generated edit script applying the Cactuar_Profits.xlsx diff
# Updates static numeric cell values (no formulas present) across all 8 sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1209.0667
$ws.Range("J19").Value = 806.36365
$ws.Range("L19").Value = 806.36365
$ws.Range("N19").Value = -1156.36365
$ws.Range("H51").Value = 6392.524
$ws.Range("J51").Value = 8103.909
$ws.Range("L51").Value = 8103.909
$ws.Range("N51").Value = -9071.909
$ws.Range("H106").Value = 37039772
$ws.Range("I106").Value = 47620420
$ws.Range("K106").Value = 47620420
$ws.Range("M106").Value = -47619789
$ws.Range("H125").Value = 4761.3335
$ws.Range("J125").Value = 4761.3335
$ws.Range("L125").Value = 42852.0015
$ws.Range("N125").Value = -47772.0015
$ws.Range("H137").Value = 6063542
$ws.Range("I137").Value = 1182.8276
$ws.Range("J137").Value = 12825405
$ws.Range("K137").Value = 3548.4828
$ws.Range("L137").Value = 38476215
$ws.Range("M137").Value = -998.4828000000002
$ws.Range("N137").Value = -38481315
$ws.Range("H138").Value = 7476.1665
$ws.Range("I138").Value = 2082.75
$ws.Range("K138").Value = 6248.25
$ws.Range("M138").Value = -1108.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3786.4753
$ws.Range("I32").Value = 2109.56
$ws.Range("K32").Value = 2109.56
$ws.Range("M32").Value = -1822.56
$ws.Range("H74").Value = 1275.081
$ws.Range("I74").Value = 1040.16
$ws.Range("J74").Value = 1764.5
$ws.Range("K74").Value = 1040.16
$ws.Range("L74").Value = 1764.5
$ws.Range("M74").Value = -166.1600000000001
$ws.Range("N74").Value = -3512.5
$ws.Range("H77").Value = 1275.081
$ws.Range("I77").Value = 1040.16
$ws.Range("J77").Value = 1764.5
$ws.Range("K77").Value = 5200.8
$ws.Range("L77").Value = 8822.5
$ws.Range("M77").Value = -832.8000000000002
$ws.Range("N77").Value = -17558.5
$ws.Range("H132").Value = 17800.35
$ws.Range("I132").Value = 18177.059
$ws.Range("K132").Value = 54531.177
$ws.Range("M132").Value = -52001.177

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 298.3846
$ws.Range("J80").Value = 214.77777
$ws.Range("L80").Value = 214.77777
$ws.Range("N80").Value = -2210.77777
$ws.Range("H83").Value = 298.3846
$ws.Range("J83").Value = 214.77777
$ws.Range("L83").Value = 1073.88885
$ws.Range("N83").Value = -11057.88885
$ws.Range("H86").Value = 1449.1111
$ws.Range("I86").Value = 1391.5
$ws.Range("J86").Value = 1495.2
$ws.Range("K86").Value = 1391.5
$ws.Range("L86").Value = 1495.2
$ws.Range("M86").Value = -268.5
$ws.Range("N86").Value = -3741.2
$ws.Range("H88").Value = 16896
$ws.Range("I88").Value = 8000
$ws.Range("K88").Value = 8000
$ws.Range("M88").Value = -7594
$ws.Range("H89").Value = 1449.1111
$ws.Range("I89").Value = 1391.5
$ws.Range("J89").Value = 1495.2
$ws.Range("K89").Value = 6957.5
$ws.Range("L89").Value = 7476
$ws.Range("M89").Value = -1341.5
$ws.Range("N89").Value = -18708
$ws.Range("H91").Value = 16896
$ws.Range("I91").Value = 8000
$ws.Range("K91").Value = 8000
$ws.Range("M91").Value = -6596
$ws.Range("H94").Value = 653230.4
$ws.Range("I94").Value = 1713048.6
$ws.Range("K94").Value = 1713048.6
$ws.Range("M94").Value = -1712597.6
$ws.Range("H107").Value = 1317.8334
$ws.Range("I107").Value = 1531
$ws.Range("J107").Value = 1230.0588
$ws.Range("K107").Value = 1531
$ws.Range("L107").Value = 1230.0588
$ws.Range("M107").Value = 389
$ws.Range("N107").Value = -5070.0588
$ws.Range("H134").Value = 2462.5574
$ws.Range("I134").Value = 2119.48
$ws.Range("K134").Value = 6358.440000000001
$ws.Range("M134").Value = -3823.440000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1279.15
$ws.Range("I16").Value = 1216.0714
$ws.Range("K16").Value = 1216.0714
$ws.Range("M16").Value = -929.0714
$ws.Range("H31").Value = 1867.4572
$ws.Range("I31").Value = 1072.3448
$ws.Range("K31").Value = 1072.3448
$ws.Range("M31").Value = -777.3448000000001
$ws.Range("H32").Value = 3502.75
$ws.Range("I32").Value = 3500
$ws.Range("K32").Value = 3500
$ws.Range("M32").Value = -3184
$ws.Range("H34").Value = 1867.4572
$ws.Range("I34").Value = 1072.3448
$ws.Range("K34").Value = 1072.3448
$ws.Range("M34").Value = -870.3448000000001
$ws.Range("H58").Value = 418543.6
$ws.Range("I58").Value = 557361.1
$ws.Range("J58").Value = 2091
$ws.Range("K58").Value = 557361.1
$ws.Range("L58").Value = 2091
$ws.Range("M58").Value = -557158.1
$ws.Range("N58").Value = -2497
$ws.Range("H113").Value = 1279.15
$ws.Range("I113").Value = 1216.0714
$ws.Range("K113").Value = 1216.0714
$ws.Range("M113").Value = 953.9286
$ws.Range("H134").Value = 2424.45
$ws.Range("I134").Value = 2443.9211
$ws.Range("K134").Value = 7331.763300000001
$ws.Range("M134").Value = -4796.763300000001
$ws.Range("H136").Value = 418543.6
$ws.Range("I136").Value = 557361.1
$ws.Range("J136").Value = 2091
$ws.Range("K136").Value = 1672083.3
$ws.Range("L136").Value = 6273
$ws.Range("M136").Value = -1669533.3
$ws.Range("N136").Value = -11373

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 429.41177
$ws.Range("I9").Value = 536.3
$ws.Range("J9").Value = 276.7143
$ws.Range("K9").Value = 1608.9
$ws.Range("L9").Value = 830.1428999999999
$ws.Range("M9").Value = -1384.9
$ws.Range("N9").Value = -1278.1429
$ws.Range("H80").Value = 5873.75
$ws.Range("J80").Value = 5873.75
$ws.Range("L80").Value = 17621.25
$ws.Range("N80").Value = -19493.25
$ws.Range("H83").Value = 5873.75
$ws.Range("J83").Value = 5873.75
$ws.Range("L83").Value = 52863.75
$ws.Range("N83").Value = -62223.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 8067.273
$ws.Range("I102").Value = 7616.5293
$ws.Range("J102").Value = 9599.799999999999
$ws.Range("K102").Value = 7616.5293
$ws.Range("L102").Value = 9599.799999999999
$ws.Range("M102").Value = -5994.5293
$ws.Range("N102").Value = -12843.8
$ws.Range("H107").Value = 3663608.5
$ws.Range("I107").Value = 6803164
$ws.Range("J107").Value = 793.8333
$ws.Range("K107").Value = 6803164
$ws.Range("L107").Value = 793.8333
$ws.Range("M107").Value = -6801244
$ws.Range("N107").Value = -4633.8333
$ws.Range("H122").Value = 552885.5
$ws.Range("I122").Value = 581868.9399999999
$ws.Range("K122").Value = 1745606.82
$ws.Range("M122").Value = -1743156.82
$ws.Range("H126").Value = 4561.4614
$ws.Range("I126").Value = 2992.2
$ws.Range("K126").Value = 8976.599999999999
$ws.Range("M126").Value = -6506.599999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1904.6842
$ws.Range("I22").Value = 932.4286
$ws.Range("J22").Value = 2471.8333
$ws.Range("K22").Value = 932.4286
$ws.Range("L22").Value = 2471.8333
$ws.Range("M22").Value = -637.4286
$ws.Range("N22").Value = -3061.8333
$ws.Range("H27").Value = 1904.6842
$ws.Range("I27").Value = 932.4286
$ws.Range("J27").Value = 2471.8333
$ws.Range("K27").Value = 932.4286
$ws.Range("L27").Value = 2471.8333
$ws.Range("M27").Value = -825.4286
$ws.Range("N27").Value = -2685.8333
$ws.Range("H136").Value = 4608.4614
$ws.Range("I136").Value = 3000.8572
$ws.Range("K136").Value = 9002.571599999999
$ws.Range("M136").Value = -6452.571599999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 13999.5
$ws.Range("I6").Value = 13999.5
$ws.Range("K6").Value = 13999.5
$ws.Range("M6").Value = -13884.5
$ws.Range("H81").Value = 1900534.8
$ws.Range("I81").Value = 1494594
$ws.Range("J81").Value = 2610931
$ws.Range("K81").Value = 2989188
$ws.Range("L81").Value = 5221862
$ws.Range("M81").Value = -2988127
$ws.Range("N81").Value = -5223984
$ws.Range("H84").Value = 1900534.8
$ws.Range("I84").Value = 1494594
$ws.Range("J84").Value = 2610931
$ws.Range("K84").Value = 14945940
$ws.Range("L84").Value = 26109310
$ws.Range("M84").Value = -14940636
$ws.Range("N84").Value = -26119918
$ws.Range("H132").Value = 15875613
$ws.Range("I132").Value = 2417513.2
$ws.Range("K132").Value = 7252539.600000001
$ws.Range("M132").Value = -7250009.600000001
$ws.Range("H136").Value = 4308.7666
$ws.Range("I136").Value = 3914.3
$ws.Range("K136").Value = 11742.9
$ws.Range("M136").Value = -9192.900000000001
